# Added better error plot formatting
# Move the "pop=50 / F = 0.9 / Cr = 0.9 / Jr=0.3" parameter labels that were
# previously duplicated on the header rows (row 1 and row 12) of each table
# down into a dedicated row (row 11) between the two tables.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: remove the parameter-label cells from the header row ---
$ws.Range("C1").ClearContents()
$ws.Range("D1").ClearContents()
$ws.Range("E1").ClearContents()
$ws.Range("M1").ClearContents()
$ws.Range("N1").ClearContents()
$ws.Range("O1").ClearContents()
$ws.Range("P1").ClearContents()

# --- Row 11: add the parameter-label cells (now on their own row) ---
$ws.Range("B11").Value = "pop=50"
$ws.Range("C11").Value = "F = 0.9"
$ws.Range("D11").Value = "Cr = 0.9"

$ws.Range("M11").Value = "pop=50"
$ws.Range("N11").Value = "F = 0.9"
$ws.Range("O11").Value = "Cr = 0.9"
$ws.Range("P11").Value = "Jr=0.3"

# --- Row 12: remove the parameter-label cells from this header row too ---
$ws.Range("C12").ClearContents()
$ws.Range("D12").ClearContents()
$ws.Range("E12").ClearContents()
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()
$ws.Range("O12").ClearContents()
$ws.Range("P12").ClearContents()

# --- Update the active selection to reflect the new layout ---
$ws.Range("K13:Q19").Select()
